$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New weekly column (18_05_2021), added one column to the right of the
# previous last column (AB -> AC), mirroring the prior week's layout.
$ws.Range("AC1").Value = "18_05_2021"

$ws.Range("AC2").Value = 205
$ws.Range("AC3").Value = 209
$ws.Range("AC4").Value = 692
$ws.Range("AC5").Value = 1050
$ws.Range("AC6").Value = 1512
$ws.Range("AC7").Value = 2369
$ws.Range("AC8").Value = 2365
$ws.Range("AC9").Value = 3228
$ws.Range("AC10").Value = 2561
$ws.Range("AC11").Value = 703

# Row 12 totals the column above it, same as every other column.
$ws.Range("AC12").Formula = "=SUM(AC2:AC11)"

# Move the selection to the newly added column, as a user would after
# typing the week's data in (mirrors the AB14 -> AC14 selection shift).
$ws.Range("AC14").Select()
